# Updated cryptos list with latest price and volume(1h) data
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "27.782.92"
$ws.Range("E2").NumberFormat = "@"
$ws.Range("E2").Value = "  +0.66%  "

$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "1.856.55"
$ws.Range("E3").NumberFormat = "@"
$ws.Range("E3").Value = "  +0.39%  "

$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "1.035"
$ws.Range("E4").NumberFormat = "@"
$ws.Range("E4").Value = "  +0.45%  "

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "322.97"
$ws.Range("E5").NumberFormat = "@"
$ws.Range("E5").Value = "  +0.66%  "

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "1.031"
$ws.Range("E6").NumberFormat = "@"
$ws.Range("E6").Value = "  +0.36%  "

$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.4400"
$ws.Range("E7").NumberFormat = "@"
$ws.Range("E7").Value = "  +0.58%  "

$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.3815"
$ws.Range("E8").NumberFormat = "@"
$ws.Range("E8").Value = "  +1.98%  "

$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.07418"
$ws.Range("E9").NumberFormat = "@"
$ws.Range("E9").Value = "  +0.20%  "

$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.8870"
$ws.Range("E10").NumberFormat = "@"
$ws.Range("E10").Value = "  +1.22%  "

$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "21.61"
$ws.Range("E11").NumberFormat = "@"
$ws.Range("E11").Value = "  +0.65%  "

$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "1.859.47"
$ws.Range("E12").NumberFormat = "@"
$ws.Range("E12").Value = "  -0.38%  "

$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "5.518"
$ws.Range("E13").NumberFormat = "@"
$ws.Range("E13").Value = "  +0.32%  "

$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "6.729"
$ws.Range("E14").NumberFormat = "@"
$ws.Range("E14").Value = "  +0.69%  "

$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "0.07193"
$ws.Range("E15").NumberFormat = "@"
$ws.Range("E15").Value = "  +0.42%  "

$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "85.38"
$ws.Range("E16").NumberFormat = "@"
$ws.Range("E16").Value = "  +3.09%  "

$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "1.038"
$ws.Range("E17").NumberFormat = "@"
$ws.Range("E17").Value = "  +0.57%  "

$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "0.000009072"
$ws.Range("E18").NumberFormat = "@"
$ws.Range("E18").Value = "  +0.44%  "

$ws.Range("E19").NumberFormat = "@"
$ws.Range("E19").Value = "  +0.49%  "

$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "15.52"
$ws.Range("E20").NumberFormat = "@"
$ws.Range("E20").Value = "  +0.63%  "

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "27.781.47"
$ws.Range("E21").NumberFormat = "@"
$ws.Range("E21").Value = "  +0.61%  "

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "5.282"
$ws.Range("E22").NumberFormat = "@"
$ws.Range("E22").Value = "  +0.58%  "

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "11.27"
$ws.Range("E23").NumberFormat = "@"
$ws.Range("E23").Value = "  +0.52%  "

$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "2.086.38"
$ws.Range("E24").NumberFormat = "@"
$ws.Range("E24").Value = "  +0.92%  "

$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "2.057"
$ws.Range("E25").NumberFormat = "@"
$ws.Range("E25").Value = "  +6.46%  "

$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "159.16"
$ws.Range("E26").NumberFormat = "@"
$ws.Range("E26").Value = "  +1.26%  "

$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "18.74"
$ws.Range("E27").NumberFormat = "@"
$ws.Range("E27").Value = "  +0.07%  "

$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "2.008"
$ws.Range("E28").NumberFormat = "@"
$ws.Range("E28").Value = "  +2.87%  "

$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "5.356"
$ws.Range("E29").NumberFormat = "@"
$ws.Range("E29").Value = "  +1.23%  "

$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "118.18"
$ws.Range("E30").NumberFormat = "@"
$ws.Range("E30").Value = "  +1.73%  "

$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "0.09095"
$ws.Range("E31").NumberFormat = "@"
$ws.Range("E31").Value = "  +0.24%  "

$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "1.214"
$ws.Range("E32").NumberFormat = "@"
$ws.Range("E32").Value = "  +0.45%  "

$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "0.7737"
$ws.Range("E33").NumberFormat = "@"
$ws.Range("E33").Value = "  +0.79%  "

$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "3.014"
$ws.Range("E34").NumberFormat = "@"
$ws.Range("E34").Value = "  +4.78%  "

$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "4.591"
$ws.Range("E35").NumberFormat = "@"
$ws.Range("E35").Value = "  +1.78%  "

$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "1.033"
$ws.Range("E36").NumberFormat = "@"
$ws.Range("E36").Value = "  +0.54%  "

$ws.Range("E37").NumberFormat = "@"
$ws.Range("E37").Value = "  +0.24%  "

$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "0.01983"
$ws.Range("E38").NumberFormat = "@"
$ws.Range("E38").Value = "  +0.13%  "

$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "0.05292"
$ws.Range("E39").NumberFormat = "@"
$ws.Range("E39").Value = "  +0.48%  "

$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "2.864"
$ws.Range("E40").NumberFormat = "@"
$ws.Range("E40").Value = "  +1.87%  "

$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "0.5196"
$ws.Range("E41").NumberFormat = "@"
$ws.Range("E41").Value = "  +0.47%  "

$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "6.902"
$ws.Range("E42").NumberFormat = "@"
$ws.Range("E42").Value = "  +2.93%  "

$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "0.1673"
$ws.Range("E43").NumberFormat = "@"
$ws.Range("E43").Value = "  -0.01%  "

$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "8.755"
$ws.Range("E44").NumberFormat = "@"
$ws.Range("E44").Value = "  +2.09%  "

$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "110.72"
$ws.Range("E45").NumberFormat = "@"
$ws.Range("E45").Value = "  +1.60%  "

$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "10.73"
$ws.Range("E46").NumberFormat = "@"
$ws.Range("E46").Value = "  +1.16%  "

$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "1.035"
$ws.Range("E47").NumberFormat = "@"
$ws.Range("E47").Value = "  +0.58%  "

$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "0.06579"
$ws.Range("E48").NumberFormat = "@"
$ws.Range("E48").Value = "  +3.11%  "

$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "1.712"
$ws.Range("E49").NumberFormat = "@"
$ws.Range("E49").Value = "  -0.26%  "

$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "0.4722"
$ws.Range("E50").NumberFormat = "@"
$ws.Range("E50").Value = "  +1.43%  "

$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "1.895"
$ws.Range("E51").NumberFormat = "@"
$ws.Range("E51").Value = "  +0.33%  "
